$wb = $excel.ActiveWorkbook

# This script applies a refreshed market-price data snapshot to the
# Mandragora_Profits crafting/leve profit tables across all job sheets.
# Each row below corresponds to one leve/recipe whose current market
# price columns (H-N) were refreshed by the scheduled data-pull runner.

$ws = $wb.Worksheets.Item("ALC")
# ALC row 4
$ws.Range("H4").Value = 1134
$ws.Range("I4").Value = 1082.2222
$ws.Range("K4").Value = 1082.2222
$ws.Range("M4").Value = -968.2221999999999

# ALC row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

# ALC row 14
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

# ALC row 33
$ws.Range("H33").Value = 541.44446
$ws.Range("I33").Value = 653.1053000000001
$ws.Range("J33").Value = 276.25
$ws.Range("K33").Value = 653.1053000000001
$ws.Range("L33").Value = 276.25
$ws.Range("M33").Value = -424.1053000000001
$ws.Range("N33").Value = -734.25

# ALC row 69
$ws.Range("H69").Value = 5193
$ws.Range("I69").Value = 6986.6665
$ws.Range("J69").Value = 4424.2856
$ws.Range("K69").Value = 20959.9995
$ws.Range("L69").Value = 13272.8568
$ws.Range("M69").Value = -20085.9995
$ws.Range("N69").Value = -15020.8568

# ALC row 72
$ws.Range("H72").Value = 5193
$ws.Range("I72").Value = 6986.6665
$ws.Range("J72").Value = 4424.2856
$ws.Range("K72").Value = 62879.9985
$ws.Range("L72").Value = 39818.5704
$ws.Range("M72").Value = -58511.9985
$ws.Range("N72").Value = -48554.5704

# ALC row 129
$ws.Range("H129").Value = 1386.4231
$ws.Range("J129").Value = 1687
$ws.Range("L129").Value = 5061
$ws.Range("N129").Value = -15061

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Range("H32").Value = 7046.647
$ws.Range("I32").Value = 7727.3887
$ws.Range("J32").Value = 4420.9287
$ws.Range("K32").Value = 7727.3887
$ws.Range("L32").Value = 4420.9287
$ws.Range("M32").Value = -7440.3887
$ws.Range("N32").Value = -4994.9287

# ARM row 74
$ws.Range("H74").Value = 1490.4884
$ws.Range("I74").Value = 1361.3667
$ws.Range("J74").Value = 1788.4615
$ws.Range("K74").Value = 1361.3667
$ws.Range("L74").Value = 1788.4615
$ws.Range("M74").Value = -487.3667
$ws.Range("N74").Value = -3536.4615

# ARM row 77
$ws.Range("H77").Value = 1490.4884
$ws.Range("I77").Value = 1361.3667
$ws.Range("J77").Value = 1788.4615
$ws.Range("K77").Value = 6806.833500000001
$ws.Range("L77").Value = 8942.307499999999
$ws.Range("M77").Value = -2438.833500000001
$ws.Range("N77").Value = -17678.3075

# ARM row 132
$ws.Range("H132").Value = 6318.5713
$ws.Range("I132").Value = 2643.4285
$ws.Range("J132").Value = 9993.714
$ws.Range("K132").Value = 7930.2855
$ws.Range("L132").Value = 29981.142
$ws.Range("M132").Value = -5400.2855
$ws.Range("N132").Value = -35041.142

$ws = $wb.Worksheets.Item("BSM")
# BSM row 24
$ws.Range("H24").Value = 1298.2222
$ws.Range("I24").Value = 461
$ws.Range("J24").Value = 2972.6667
$ws.Range("K24").Value = 461
$ws.Range("L24").Value = 2972.6667
$ws.Range("M24").Value = -226
$ws.Range("N24").Value = -3442.6667

# BSM row 105
$ws.Range("H105").Value = 3023.5454
$ws.Range("I105").Value = 2919.889
$ws.Range("K105").Value = 2919.889
$ws.Range("M105").Value = -1172.889

# BSM row 134
$ws.Range("H134").Value = 7908.1924
$ws.Range("I134").Value = 4625.625
$ws.Range("J134").Value = 9367.111000000001
$ws.Range("K134").Value = 13876.875
$ws.Range("L134").Value = 28101.333
$ws.Range("M134").Value = -11341.875
$ws.Range("N134").Value = -33171.333

$ws = $wb.Worksheets.Item("CRP")
# CRP row 50
$ws.Range("H50").Value = 22000
$ws.Range("J50").Value = 22000
$ws.Range("L50").Value = 22000
$ws.Range("N50").Value = -23250

# CRP row 54
$ws.Range("H54").Value = 5000
$ws.Range("J54").Value = 5000
$ws.Range("L54").Value = 5000
$ws.Range("N54").Value = -6316

# CRP row 59
$ws.Range("H59").Value = 32873.223
$ws.Range("I59").Value = 15000
$ws.Range("J59").Value = 35107.375
$ws.Range("K59").Value = 15000
$ws.Range("L59").Value = 35107.375
$ws.Range("M59").Value = -13855
$ws.Range("N59").Value = -37397.375

$ws = $wb.Worksheets.Item("CUL")
# CUL row 69
$ws.Range("H69").Value = 2778
$ws.Range("J69").Value = 3272.7273
$ws.Range("L69").Value = 9818.1819
$ws.Range("N69").Value = -11440.1819

# CUL row 72
$ws.Range("H72").Value = 2778
$ws.Range("J72").Value = 3272.7273
$ws.Range("L72").Value = 29454.5457
$ws.Range("N72").Value = -37566.5457

$ws = $wb.Worksheets.Item("GSM")
# GSM row 18
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# GSM row 31
$ws.Range("H31").Value = 755.125
$ws.Range("I31").Value = 755.125
$ws.Range("K31").Value = 755.125
$ws.Range("M31").Value = -463.125

# GSM row 37
$ws.Range("H37").Value = 755.125
$ws.Range("I37").Value = 755.125
$ws.Range("K37").Value = 755.125
$ws.Range("M37").Value = -478.125

# GSM row 43
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

# GSM row 46
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

# GSM row 57
$ws.Range("H57").Value = 9666.666999999999
$ws.Range("I57").Value = 9666.666999999999
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 9666.666999999999
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -8846.666999999999
$ws.Range("N57").ClearContents()

# GSM row 80
$ws.Range("H80").Value = 2635.625
$ws.Range("I80").Value = 2635.625
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2635.625
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1637.625
$ws.Range("N80").ClearContents()

# GSM row 83
$ws.Range("H83").Value = 2635.625
$ws.Range("I83").Value = 2635.625
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 13178.125
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -8186.125
$ws.Range("N83").Value = 0

# GSM row 132
$ws.Range("H132").Value = 3475372.8
$ws.Range("I132").Value = 10420044
$ws.Range("J132").Value = 3037.25
$ws.Range("K132").Value = 31260132
$ws.Range("L132").Value = 9111.75
$ws.Range("M132").Value = -31257602
$ws.Range("N132").Value = -14171.75

$ws = $wb.Worksheets.Item("LTW")
# LTW row 16
$ws.Range("H16").Value = 2464.4092
$ws.Range("I16").Value = 2579.6428
$ws.Range("J16").Value = 2262.75
$ws.Range("K16").Value = 2579.6428
$ws.Range("L16").Value = 2262.75
$ws.Range("M16").Value = -2409.6428
$ws.Range("N16").Value = -2602.75

# LTW row 40
$ws.Range("H40").Value = 111119090
$ws.Range("I40").Value = 250003200
$ws.Range("J40").Value = 11798
$ws.Range("K40").Value = 250003200
$ws.Range("L40").Value = 11798
$ws.Range("M40").Value = -250003064
$ws.Range("N40").Value = -12070

# LTW row 46
$ws.Range("H46").Value = 1429262.9
$ws.Range("I46").Value = 570
$ws.Range("J46").Value = 5000995
$ws.Range("K46").Value = 570
$ws.Range("L46").Value = 5000995
$ws.Range("M46").Value = -382
$ws.Range("N46").Value = -5001371

# LTW row 74
$ws.Range("H74").Value = 16879.572
$ws.Range("I74").Value = 12725.667
$ws.Range("J74").Value = 19995
$ws.Range("K74").Value = 12725.667
$ws.Range("L74").Value = 19995
$ws.Range("M74").Value = -11727.667
$ws.Range("N74").Value = -21991

# LTW row 77
$ws.Range("H77").Value = 16879.572
$ws.Range("I77").Value = 12725.667
$ws.Range("J77").Value = 19995
$ws.Range("K77").Value = 38177.001
$ws.Range("L77").Value = 59985
$ws.Range("M77").Value = -33185.001
$ws.Range("N77").Value = -69969

$ws = $wb.Worksheets.Item("WVR")
# WVR row 100
$ws.Range("H100").Value = 685.1852
$ws.Range("I100").Value = 443.63635
$ws.Range("J100").Value = 1748
$ws.Range("K100").Value = 887.2727
$ws.Range("L100").Value = 3496
$ws.Range("M100").Value = -346.2727
$ws.Range("N100").Value = -4578

# WVR row 107
$ws.Range("H107").Value = 1021.1667
$ws.Range("I107").Value = 781
$ws.Range("J107").Value = 2222
$ws.Range("K107").Value = 2343
$ws.Range("L107").Value = 6666
$ws.Range("M107").Value = -423
$ws.Range("N107").Value = -10506

# WVR row 136
$ws.Range("H136").Value = 6411933
$ws.Range("I136").Value = 10417774
$ws.Range("J136").Value = 2586.6667
$ws.Range("K136").Value = 31253322
$ws.Range("L136").Value = 7760.000100000001
$ws.Range("M136").Value = -31250772
$ws.Range("N136").Value = -12860.0001
